# dmb ran run5 samples for March 2020, ran first 9, 3 currently running
#
# Append the new CRM-accuracy data row (row 35) for sample date 2021-05-11,
# following the same layout as the existing rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = 20210511
$ws.Range("B35").Value = 2228.8789999999999
$ws.Range("C35").Value = 2224.4699999999998
$ws.Range("D35").Formula = "=100*(B35-C35)/C35"
$ws.Range("E35").Value = 180
$ws.Range("F35").Value = "CRM opened 20210418"

# Leave the selection where data entry finished, matching the author's
# cursor position after typing the new row.
$ws.Range("D35").Select()
